# Regenerate the "K" column (column G) values for every data row on
# Sheet1. The source data (Strike# -> K) was recomputed upstream; this
# script writes the newly computed K values into the existing sheet,
# leaving every other column untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New K (column G) values, in row order starting at sheet row 2
# (row index 0 in the original "save_data" table) through row 59
# (row index 57).
$kValues = @(
    0,1,2,1,1,1,4,0,1,0,
    1,1,0,2,0,1,1,1,2,2,
    0,1,0,3,1,3,1,0,2,4,
    1,1,2,1,2,0,1,1,1,1,
    1,1,3,1,0,3,2,4,2,5,
    4,1,1,7,3,1,1,1
)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
